$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 169, shifting existing rows 169:247 down to 170:248
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new weekly record
$ws.Cells.Item(169, 1).Value = 10
$ws.Cells.Item(169, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(169, 3).Value = "La Araucanía"
$ws.Cells.Item(169, 4).Value = 44704
$ws.Cells.Item(169, 5).Value = 9
$ws.Cells.Item(169, 6).Value = 100112039
$ws.Cells.Item(169, 7).Value = "Ciboulette"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 30
$ws.Cells.Item(169, 11).Value = 6000
$ws.Cells.Item(169, 12).Value = 6000
$ws.Cells.Item(169, 13).Value = 6000
$ws.Cells.Item(169, 14).Value = "$/docena de atados"
$ws.Cells.Item(169, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(169, 16).Value = 2000
$ws.Cells.Item(169, 17).Value = 3
$ws.Cells.Item(169, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Cells.Item(169, 4).NumberFormat = $ws.Cells.Item(170, 4).NumberFormat
